# Fixed update to excel issue
# - Rename the "Requested quantity" headers on the existing sheets
# - Add a new "PO Forecast" sheet (after "Monthly Trend") with forecast data

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Rename the "Requested quantity" column headers
$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet right after "Monthly Trend"
$newWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$newWs.Name = "PO Forecast"

# Match page margins used by the rest of the workbook
$newWs.PageSetup.LeftMargin = $ws1.PageSetup.LeftMargin
$newWs.PageSetup.RightMargin = $ws1.PageSetup.RightMargin
$newWs.PageSetup.TopMargin = $ws1.PageSetup.TopMargin
$newWs.PageSetup.BottomMargin = $ws1.PageSetup.BottomMargin
$newWs.PageSetup.HeaderMargin = $ws1.PageSetup.HeaderMargin
$newWs.PageSetup.FooterMargin = $ws1.PageSetup.FooterMargin

# Reuse the existing header formatting (bold font, border, centered) for the
# new sheet's header row, and the existing date-column formatting for column A
$ws1.Range("A1:B1").Copy()
$newWs.Range("A1:B1").PasteSpecial(-4122)
$newWs.Range("C1:D1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$newWs.Range("A2:A16").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Header row
$newWs.Range("A1").Value = "ds"
$newWs.Range("B1").Value = "PO_Forecast"
$newWs.Range("C1").Value = "yhat_lower"
$newWs.Range("D1").Value = "yhat_upper"

# Forecast data rows
$data = @(
    @(44934.99999999999, 20, 6.981058702900047, 32.11578742938477),
    @(44955.99999999999, 20, 7.702238846755778, 34.06030642138945),
    @(44976.99999999999, 21, 8.112867433642554, 32.87319517884598),
    @(44983.99999999999, 21, 8.377510236890767, 33.75407578564982),
    @(45011.99999999999, 21, 7.67581306409935, 32.81615794427092),
    @(45060.99999999999, 22, 8.094381295522812, 34.52281938476338),
    @(45319.99999999999, 24, 10.17671886918076, 36.72900428442777),
    @(45326.99999999999, 24, 11.52024024019917, 37.59953433385995),
    @(45333.99999999999, 25, 11.26609193016053, 37.6434247588341),
    @(45340.99999999999, 25, 11.94733405337301, 36.97762333267378),
    @(45347.99999999999, 25, 11.98094862188953, 37.74527960051481),
    @(45354.99999999999, 25, 11.56551361149244, 37.76184144009169),
    @(45361.99999999999, 25, 11.70195170217422, 38.41510869714959),
    @(45368.99999999999, 25, 11.44424078725643, 38.60338586921919),
    @(45375.99999999999, 25, 11.9804654564349, 37.72873540014)
)

$row = 2
foreach ($r in $data) {
    $newWs.Cells.Item($row, 1).Value = $r[0]
    $newWs.Cells.Item($row, 2).Value = $r[1]
    $newWs.Cells.Item($row, 3).Value = $r[2]
    $newWs.Cells.Item($row, 4).Value = $r[3]
    $row++
}

Write-Output "PO Forecast sheet added with $($data.Count) data rows"
